# Normalize ObjTables header attribute names to lowerCamelCase
# (Type -> type, Id -> id, ObjTablesVersion -> objTablesVersion)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("!!Test")
$ws.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws.Range("A2").Value = "!!ObjTables type='Data' id='Test'"

$ws = $wb.Worksheets.Item("!!Deleted models")
$ws.Range("A1").Value = "!!ObjTables type='Data' id='DeletedModel'"

$ws = $wb.Worksheets.Item("!!Property")
$ws.Range("A1").Value = "!!ObjTables type='Data' id='Property'"

$ws = $wb.Worksheets.Item("!!Subtests")
$ws.Range("A1").Value = "!!ObjTables type='Data' id='Subtest'"

$ws = $wb.Worksheets.Item("!!References")
$ws.Range("A1").Value = "!!ObjTables type='Data' id='Reference'"
